$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet right before "2022-Q3" so the tab
#    order becomes: 总计, 2022-Q4, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4,
#    2021-Q3, 2021-Q2, 2021-Q1
# ---------------------------------------------------------------------------
$sheets = $wb.Worksheets
$q3 = $sheets.Item("2022-Q3")
$q4 = $sheets.Add($q3)
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# 2) Borrow the exact cell-formatting skeleton (bold/border header row +
#    bold/border index column) used by the other quarter sheets. "2022-Q1"
#    has 8 data rows, i.e. exactly A1:H8, the same footprint the new sheet
#    needs.
# ---------------------------------------------------------------------------
$q1 = $sheets.Item("2022-Q1")
$q1.Range("A1:H8").Copy()
$q4.Range("A1:H8").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Header row
# ---------------------------------------------------------------------------
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------------
# 4) Data rows. Columns B, D, E, F, G hold text that merely looks numeric
#    (fund codes/leading zeros, percentages kept as strings in the source
#    data) so force text storage via NumberFormat "@" before assigning.
#    Columns A (row index) and H (rank) are genuine numbers.
# ---------------------------------------------------------------------------
$codes   = @("001743","010296","590006","010297","003456","002409","002410")
$names   = @("诺安优选回报灵活配置混合","万家互联互通中国优势量化策略混合A","中邮中小盘灵活配置混合","万家互联互通中国优势量化策略混合C","信澳新目标灵活配置混合","华夏新活力灵活配置混合A","华夏新活力灵活配置混合C")
$scale   = @("21.48","4.37","2.56","0.47","0.44","0.12","0.00")
$stock   = @("72.70","94.52","74.40","94.52","51.24","77.58","77.58")
$pct     = @("3.16","5.91","2.47","5.91","1.09","3.25","3.25")
$heldTxt = @("0.6788","0.2583","0.0632","0.0278","0.0048","0.0039",$null)
$rank    = @(4,2,3,2,7,10,10)

$textCols = @("B","D","E","F","G")
for ($i = 0; $i -lt 7; $i++) {
    $r = $i + 2
    foreach ($col in $textCols) {
        $q4.Range("$col$r").NumberFormat = "@"
    }

    $q4.Range("A$r").Value = $i
    $q4.Range("B$r").Value = $codes[$i]
    $q4.Range("C$r").Value = $names[$i]
    $q4.Range("D$r").Value = $scale[$i]
    $q4.Range("E$r").Value = $stock[$i]
    $q4.Range("F$r").Value = $pct[$i]
    if ($heldTxt[$i] -eq $null) {
        # row 8 (华夏新活力灵活配置混合C) stores the held-value as a real 0,
        # not as text, unlike every other row in this column.
        $q4.Range("G$r").NumberFormat = "General"
        $q4.Range("G$r").Value = 0
    } else {
        $q4.Range("G$r").Value = $heldTxt[$i]
    }
    $q4.Range("H$r").Value = $rank[$i]
}

# ---------------------------------------------------------------------------
# 5) "总计" summary sheet: insert a new row 2 for the 2022-Q4 totals and
#    push the existing rows down (their own values are untouched).
# ---------------------------------------------------------------------------
$zj = $sheets.Item("总计")
$zj.Rows.Item(2).Insert()
$zj.Range("B2:D2").ClearFormats()

# Re-apply the row-index column styling (bold/border) that column A carries
# on every other row of this sheet.
$zj.Range("A3").Copy()
$zj.Range("A2").PasteSpecial(-4122)

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q4"
$zj.Range("C2").Value = 7
$zj.Range("D2").Value = 1.04
